$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record for Puerro (Vega Modelo de Temuco) was inserted
# at row 277, pushing all the following records (old 277..332) down by one
# row (they now occupy 278..333). This mirrors the commit "Fruta / hortaliza,
# semanal" which prepends the latest weekly observation to the dataset.

$ws.Range("A277:R277").Insert()

$ws.Range("A277").Value = 10
$ws.Range("B277").Value = "Vega Modelo de Temuco"
$ws.Range("C277").Value = "La Araucanía"
$ws.Range("D277").Value = 45173
$ws.Range("E277").Value = 9
$ws.Range("F277").Value = 100112005
$ws.Range("G277").Value = "Puerro"
$ws.Range("H277").Value = "Azul de Maquehue"
$ws.Range("I277").Value = "Primera"
$ws.Range("J277").Value = 100
$ws.Range("K277").Value = 9000
$ws.Range("L277").Value = 9000
$ws.Range("M277").Value = 9000
$ws.Range("N277").Value = "$/docena de paquetes"
$ws.Range("O277").Value = "Provincia de Cautín"
$ws.Range("P277").Value = 750
$ws.Range("Q277").Value = 12
$ws.Range("R277").Value = "Hortaliza"
